$d = $word.ActiveDocument

# --- Programa paragraph: insert a manual line break (<w:br/>) at each bullet boundary ---
$d.Content.Find.Execute(". Cidades como sistemas socioecológicos;- Ambiente físi", $true, $false, $false, $false, $false, $true, 1, $false, ". Cidades como sistemas socioecológicos;^l- Ambiente físi", 2) | Out-Null
$d.Content.Find.Execute("o da vegetação na atenuação da poluição;- O metabolismo", $true, $false, $false, $false, $false, $true, 1, $false, "o da vegetação na atenuação da poluição;^l- O metabolismo", 2) | Out-Null
$d.Content.Find.Execute(" metabolismo urbano (Abel Woman, 1965); - Biodiversidad", $true, $false, $false, $false, $false, $true, 1, $false, " metabolismo urbano (Abel Woman, 1965); ^l- Biodiversidad", 2) | Out-Null
$d.Content.Find.Execute("ersidade urbana às mudanças climáticas);- Infraestrutur", $true, $false, $false, $false, $false, $true, 1, $false, "ersidade urbana às mudanças climáticas);^l- Infraestrutur", 2) | Out-Null
$d.Content.Find.Execute("borização urbana em cidades brasileiras;- Infraestrutur", $true, $false, $false, $false, $false, $true, 1, $false, "borização urbana em cidades brasileiras;^l- Infraestrutur", 2) | Out-Null
$d.Content.Find.Execute("e biomonitoramento da qualidade da água;- Heterogeneida", $true, $false, $false, $false, $false, $true, 1, $false, "e biomonitoramento da qualidade da água;^l- Heterogeneida", 2) | Out-Null
$d.Content.Find.Execute("cidades aos eventos climáticos extremos;- Agenda 2030 e", $true, $false, $false, $false, $false, $true, 1, $false, "cidades aos eventos climáticos extremos;^l- Agenda 2030 e", 2) | Out-Null
$d.Content.Find.Execute(" para os ODS 3, ODS 6, ODS 11 e ODS 13; - Saída para ca", $true, $false, $false, $false, $false, $true, 1, $false, " para os ODS 3, ODS 6, ODS 11 e ODS 13; ^l- Saída para ca", 2) | Out-Null
$d.Content.Find.Execute("variáveis socioeconômicas e ambientais);- Apresentação ", $true, $false, $false, $false, $false, $true, 1, $false, "variáveis socioeconômicas e ambientais);^l- Apresentação ", 2) | Out-Null

# --- Bibliografia paragraph: insert a blank line (two manual line breaks) between citations ---
$d.Content.Find.Execute("lo, SP: Oficina de Textos, 384 p., 2015.Elmqvist T, Fra", $true, $false, $false, $false, $false, $true, 1, $false, "lo, SP: Oficina de Textos, 384 p., 2015.^l^lElmqvist T, Fra", 2) | Out-Null
$d.Content.Find.Execute("sessment. Springer Nature, 755 p., 2013.Gaston KJ. Urba", $true, $false, $false, $false, $false, $true, 1, $false, "sessment. Springer Nature, 755 p., 2013.^l^lGaston KJ. Urba", 2) | Out-Null
$d.Content.Find.Execute("y; New York: Cambridge University Press.James P, Dougla", $true, $false, $false, $false, $false, $true, 1, $false, "y; New York: Cambridge University Press.^l^lJames P, Dougla", 2) | Out-Null
$d.Content.Find.Execute("ion. 2ª Edição. Routledge, 530 p., 2023.Niemelä J, Breu", $true, $false, $false, $false, $false, $true, 1, $false, "ion. 2ª Edição. Routledge, 530 p., 2023.^l^lNiemelä J, Breu", 2) | Out-Null
$d.Content.Find.Execute(" applications. OUP Oxford, 859 p., 2011.Ricklefs RE. A ", $true, $false, $false, $false, $false, $true, 1, $false, " applications. OUP Oxford, 859 p., 2011.^l^lRicklefs RE. A ", 2) | Out-Null
$d.Content.Find.Execute("ogan. Rio de Janeiro, RJ, 1664 p., 2021.Townsend CR, Be", $true, $false, $false, $false, $false, $true, 1, $false, "ogan. Rio de Janeiro, RJ, 1664 p., 2021.^l^lTownsend CR, Be", 2) | Out-Null
$d.Content.Find.Execute("Editora. Porto Alegre, RS, 632 p., 2018.Bibliografia co", $true, $false, $false, $false, $false, $true, 1, $false, "Editora. Porto Alegre, RS, 632 p., 2018.^l^lBibliografia co", 2) | Out-Null
$d.Content.Find.Execute("Bibliografia complementar:Barbosa VL, Jún", $true, $false, $false, $false, $false, $true, 1, $false, "Bibliografia complementar:^l^lBarbosa VL, Jún", 2) | Out-Null
$d.Content.Find.Execute("tal. Geografia (Londrina), 18(2), 21-36.Douglas I (2012", $true, $false, $false, $false, $false, $true, 1, $false, "tal. Geografia (Londrina), 18(2), 21-36.^l^lDouglas I (2012", 2) | Out-Null
$d.Content.Find.Execute("ronmental Sustainability, 4(4), 385-392.Ferreira ML, Ba", $true, $false, $false, $false, $false, $true, 1, $false, "ronmental Sustainability, 4(4), 385-392.^l^lFerreira ML, Ba", 2) | Out-Null
$d.Content.Find.Execute("ling. Ecological Processes, 10(1), 1-13.Ferreira ML, De", $true, $false, $false, $false, $false, $true, 1, $false, "ling. Ecological Processes, 10(1), 1-13.^l^lFerreira ML, De", 2) | Out-Null
$d.Content.Find.Execute("em services. Sustainability, 10(3), 684.Ferreira ML, Ri", $true, $false, $false, $false, $false, $true, 1, $false, "em services. Sustainability, 10(3), 684.^l^lFerreira ML, Ri", 2) | Out-Null
